$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1. Add the new changelog entry as row 20, and an extra blank row (21)
#    below it -- mirrors the OOXML diff which grows the sheet from
#    A1:C19 to A1:C21.
# ---------------------------------------------------------------------

# Copy the formatting (styles) of the previous data row (19) down onto
# the new row 20 so the new cells pick up the same styles (left/top
# aligned text, wrapped "Details" column, date-formatted "Date" column)
# without inventing any new cellXfs entries.
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)

# The "Date" column on row 19 is blank, so it does not carry the date
# number-format style -- pull that from a row that actually has a date
# (row 17) instead, just for the C column.
$ws.Range("C17").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new changelog row.
$ws.Range("A20").Value = "[1.17.5]"
$ws.Range("B20").Value = "~ adjusting Envelop template to align the recepient name at center" + [char]10 + "~ upgrading the Spire dlls to hotfix 6.8.11"
$ws.Range("C20").Value = 43341

# Row 19 grows from a 3-line entry to a shorter 2-line one in the new
# version, so its height shrinks from 90 to 65.25.
$ws.Rows.Item(19).RowHeight = 65.25

# New row 20 is a 2-line entry, same as rows 17/18 -- 30pt tall.
$ws.Rows.Item(20).RowHeight = 30

# Row 21: a trailing blank row, styled like other blank-ish rows
# (plain left/top alignment on all three columns, e.g. row 5).
$ws.Range("A5:C5").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Grow the table ("Table2") so it covers the two new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C21"))

# ---------------------------------------------------------------------
# 3. Update the sheet view: freeze the header row, scroll down so row
#    14 is the first visible row under the freeze line, and select
#    B18:B20 (active cell B18) -- matches the refreshed view state in
#    the diff.
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.Goto($ws.Range("A14"))
$ws.Range("B18:B20").Select()
